$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.171.36"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.669.86"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "210.58"
$ws.Range("E5").Value = "  -2.56%  "
$ws.Range("D6").Value = "0.5209"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "0.2625"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "0.06322"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "21.20"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "0.07551"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "1.675.93"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "4.442"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "0.5494"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").Value = "0.000008023"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "66.38"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "26.179.34"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "4.754"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("D20").Value = "187.15"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "10.34"
$ws.Range("E21").Value = "  -4.08%  "
$ws.Range("D22").Value = "6.211"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "149.87"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "0.1241"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "7.490"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").Value = "15.82"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "0.06340"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "1.353"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").Value = "1.283"
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("D31").Value = "3.520"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "3.414"
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").Value = "1.644"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").Value = "1.005"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "0.6045"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").Value = "2.405"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "2.754"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "1.112.43"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "6.127"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "0.01614"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "0.8654"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "100.39"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "1.824.00"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "55.61"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "8.087"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "0.4241"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "5.929"
$ws.Range("E51").Value = "  -1.31%  "
